$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.200.90'
$ws.Range('E2').Value = '  -4.76%  '
$ws.Range('D3').Value = '2.534.90'
$ws.Range('E3').Value = '  -4.38%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '504.24'
$ws.Range('E5').Value = '  -5.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.89'
$ws.Range('E6').Value = '  -8.43%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.561'
$ws.Range('E8').Value = '  -5.20%  '
$ws.Range('D9').Value = '2.538.81'
$ws.Range('E9').Value = '  -4.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.08'
$ws.Range('E10').Value = '  -8.46%  '
$ws.Range('E11').Value = '  -7.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.330'
$ws.Range('E12').Value = '  -6.29%  '
$ws.Range('E13').Value = '  -0.68%  '
$ws.Range('D14').Value = '2.981.89'
$ws.Range('E14').Value = '  -4.24%  '
$ws.Range('D15').Value = '58.197.99'
$ws.Range('E15').Value = '  -4.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.61'
$ws.Range('E16').Value = '  -6.79%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000134'
$ws.Range('E17').Value = '  -7.06%  '
$ws.Range('D18').Value = '2.543.75'
$ws.Range('E18').Value = '  -3.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.49'
$ws.Range('E19').Value = '  -5.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '340.38'
$ws.Range('E20').Value = '  -4.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.01'
$ws.Range('E21').Value = '  -6.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.91'
$ws.Range('E23').Value = '  -5.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.37'
$ws.Range('E24').Value = '  -2.17%  '
$ws.Range('E25').Value = '  -5.69%  '
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('D27').Value = '2.655.53'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.157'
$ws.Range('E28').Value = '  -6.74%  '
$ws.Range('D29').Value = '0.0₃0776'
$ws.Range('E29').Value = '  -9.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.89'
$ws.Range('E30').Value = '  -6.68%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '149.61'
$ws.Range('E32').Value = '  -0.41%  '
$ws.Range('B33').Value = 'Aptos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.84'
$ws.Range('E33').Value = '  -5.95%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.45'
$ws.Range('E34').Value = '  -5.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.52'
$ws.Range('E35').Value = '  -7.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.86'
$ws.Range('E36').Value = '  -6.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.894'
$ws.Range('E37').Value = '  -1.00%  '
$ws.Range('E38').Value = '  -8.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.91'
$ws.Range('E39').Value = '  -2.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.811'
$ws.Range('E40').Value = '  -12.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.37'
$ws.Range('E41').Value = '  -8.87%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '280.26'
$ws.Range('E42').Value = '  -9.67%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.50'
$ws.Range('E43').Value = '  -8.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.999'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0988'
$ws.Range('E45').Value = '  -3.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.598'
$ws.Range('E46').Value = '  -7.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0529'
$ws.Range('E47').Value = '  -6.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.28'
$ws.Range('E48').Value = '  -0.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.51'
$ws.Range('E49').Value = '  -7.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0224'
$ws.Range('E50').Value = '  -6.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.49'
$ws.Range('E51').Value = '  -11.18%  '
